$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 192, shifting existing rows 192:209 down to 193:210
$ws.Rows.Item(192).Insert()

# Populate the new row 192 with data (same constant columns as surrounding rows,
# plus the new specific values for D, J, K, L, M, O, P)
$ws.Cells.Item(192, 1).Value2 = 4
$ws.Cells.Item(192, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(192, 3).Value2 = "Los Lagos"
$ws.Cells.Item(192, 4).Value2 = 44578
$ws.Cells.Item(192, 5).Value2 = 10
$ws.Cells.Item(192, 6).Value2 = 100112043
$ws.Cells.Item(192, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(192, 8).Value2 = "Sin especificar"
$ws.Cells.Item(192, 9).Value2 = "Primera"
$ws.Cells.Item(192, 10).Value2 = 200
$ws.Cells.Item(192, 11).Value2 = 12000
$ws.Cells.Item(192, 12).Value2 = 12000
$ws.Cells.Item(192, 13).Value2 = 12000
$ws.Cells.Item(192, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(192, 15).Value2 = "Región del Maule"
$ws.Cells.Item(192, 16).Value2 = 200
$ws.Cells.Item(192, 17).Value2 = 60
$ws.Cells.Item(192, 18).Value2 = "Hortaliza"
